$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.283.41"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.81"
$ws.Range("E3").Value = "  -2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.38"
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  -2.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2881"
$ws.Range("E8").Value = "  -3.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06599"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.889.33"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.86"
$ws.Range("E11").Value = "  -2.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07374"
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.187"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.92"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6611"
$ws.Range("E15").Value = "  -2.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.258.72"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.48"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007724"
$ws.Range("E19").Value = "  -2.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.458"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.146.30"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "192.63"
$ws.Range("E23").Value = "  -3.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.201"
$ws.Range("E24").Value = "  -2.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.434"
$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.28"
$ws.Range("E26").Value = "  +1.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.943"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.445"
$ws.Range("E29").Value = "  -2.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.268"
$ws.Range("E30").Value = "  -2.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09162"
$ws.Range("E31").Value = "  -0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05073"
$ws.Range("E33").Value = "  -4.28%  "

$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.143"
$ws.Range("E35").Value = "  +1.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.710"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01833"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.632"
$ws.Range("E38").Value = "  -3.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9153"
$ws.Range("E39").Value = "  -1.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.079"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.56"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4332"
$ws.Range("E42").Value = "  -3.34%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.884"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.716"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1358"
$ws.Range("E46").Value = "  -3.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.577"
$ws.Range("E47").Value = "  +8.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.22"
$ws.Range("E48").Value = "  -9.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.934"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.20"
$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05728"
$ws.Range("E51").Value = "  -2.62%  "
